$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Fill in the "完成情况" (progress) percentage for the week-13 block
#    (rows 233-238, column C) and switch those cells to a percent format.
# ---------------------------------------------------------------------------
$progressRows = 233,234,235,236,237,238
foreach ($r in $progressRows) {
    $c = $ws.Cells.Item($r, 3)
    $c.Value = 0.4
    $c.NumberFormat = "0%"
}

# ---------------------------------------------------------------------------
# 2) Append a new weekly block (rows 241-250) for
#    "日期：2018.12.3 第十四周周一", mirroring the existing week-13 block
#    (rows 231-240) both in data and in formatting.
# ---------------------------------------------------------------------------

# Row 241: week header (merged A:D), same look as row 231.
$ws.Range("A231").Copy() | Out-Null
$ws.Range("A241").PasteSpecial(-4122) | Out-Null
$ws.Range("B231").Copy() | Out-Null
$ws.Range("B241").PasteSpecial(-4122) | Out-Null
$ws.Range("C231").Copy() | Out-Null
$ws.Range("C241").PasteSpecial(-4122) | Out-Null
$ws.Range("D231").Copy() | Out-Null
$ws.Range("D241").PasteSpecial(-4122) | Out-Null
$ws.Range("A241").Value = "日期：2018.12.3 第十四周周一"

# Row 242: column headers, same look as row 232.
$ws.Range("A232").Copy() | Out-Null
$ws.Range("A242").PasteSpecial(-4122) | Out-Null
$ws.Range("B232").Copy() | Out-Null
$ws.Range("B242").PasteSpecial(-4122) | Out-Null
$ws.Range("C232").Copy() | Out-Null
$ws.Range("C242").PasteSpecial(-4122) | Out-Null
$ws.Range("D232").Copy() | Out-Null
$ws.Range("D242").PasteSpecial(-4122) | Out-Null
$ws.Range("A242").Value = "组员"
$ws.Range("B242").Value = "计划内容"
$ws.Range("C242").Value = "完成情况"
$ws.Range("D242").Value = "备注"

# Rows 243-248: member rows, same look/content as rows 233-238, but the
# progress column is left blank (format-only, no value).
$memberSrcRows = 233,234,235,236,237,238
$memberDstRows = 243,244,245,246,247,248
$memberNames = "练富珊","黄成志","黄皓燊","郑嘉蔚","陈碧容","辛伟达"
$memberTasks = "网络交互文档[前半部分]","网络交互文档[后半部分]","测试文档[前半部分]","测试文档[后半部分]","使用手册[前半部分]","使用手册[后半部分]"

for ($i = 0; $i -lt $memberSrcRows.Length; $i++) {
    $src = $memberSrcRows[$i]
    $dst = $memberDstRows[$i]

    $ws.Cells.Item($src, 1).Copy() | Out-Null
    $ws.Cells.Item($dst, 1).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($src, 2).Copy() | Out-Null
    $ws.Cells.Item($dst, 2).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($src, 3).Copy() | Out-Null
    $ws.Cells.Item($dst, 3).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($src, 4).Copy() | Out-Null
    $ws.Cells.Item($dst, 4).PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($dst, 3).ClearContents() | Out-Null
}

$ws.Range("A243").Value = "练富珊"
$ws.Range("A244").Value = "黄成志"
$ws.Range("A245").Value = "黄皓燊"
$ws.Range("A246").Value = "郑嘉蔚"
$ws.Range("A247").Value = "陈碧容"
$ws.Range("A248").Value = "辛伟达"

$ws.Range("B243").Value = "网络交互文档[前半部分]"
$ws.Range("B244").Value = "网络交互文档[后半部分]"
$ws.Range("B245").Value = "测试文档[前半部分]"
$ws.Range("B246").Value = "测试文档[后半部分]"
$ws.Range("B247").Value = "使用手册[前半部分]"
$ws.Range("B248").Value = "使用手册[后半部分]"

# Row 249-250: "总结：" block, same look as rows 239-240.
$ws.Range("A239").Copy() | Out-Null
$ws.Range("A249").PasteSpecial(-4122) | Out-Null
$ws.Range("B239").Copy() | Out-Null
$ws.Range("B249").PasteSpecial(-4122) | Out-Null
$ws.Range("C239").Copy() | Out-Null
$ws.Range("C249").PasteSpecial(-4122) | Out-Null
$ws.Range("D239").Copy() | Out-Null
$ws.Range("D249").PasteSpecial(-4122) | Out-Null
$ws.Range("A249").Value = "总结："

$ws.Range("A240").Copy() | Out-Null
$ws.Range("A250").PasteSpecial(-4122) | Out-Null
$ws.Range("B240").Copy() | Out-Null
$ws.Range("B250").PasteSpecial(-4122) | Out-Null
$ws.Range("C240").Copy() | Out-Null
$ws.Range("C250").PasteSpecial(-4122) | Out-Null
$ws.Range("D240").Copy() | Out-Null
$ws.Range("D250").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# 3) Merge the new header / summary rows, matching the pattern used by every
#    other weekly block.
# ---------------------------------------------------------------------------
$ws.Range("A241:D241").Merge() | Out-Null
$ws.Range("A249:D250").Merge() | Out-Null

# ---------------------------------------------------------------------------
# 4) Update the view so the newly added block is in focus.
# ---------------------------------------------------------------------------
$ws.Range("A239:D240").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 223
$excel.ActiveWindow.ScrollColumn = 1
